# Edit "TextBox 60" (Testing section, bottom-right) on slide 1:
#  - resize/reposition the box
#  - tweak bullet #2 text
#  - insert two new bullets
#  - split the final bullet into three runs

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 60") {
        $shp = $s.Shapes.Item($i)
        break
    }
}

# --- Resize / reposition (EMU 26169410,21642497 9673165x5632311 expressed in points) ---
$shp.Left = 2060.5834646
$shp.Top = 1704.133622
$shp.Width = 761.6665354
$shp.Height = 443.4890601

$tr = $shp.TextFrame.TextRange

# --- Bullet 2: tweak wording (replace whole run so formatting/rPr is preserved as a single run) ---
$para2 = $tr.Paragraphs(2)
$full2 = $para2.Characters(1, $para2.Length)
$full2.Text = "No unit testing for main code because Black Box tests cover the majority of use cases"

# --- Insert both new bullets right after bullet 2 (insert the *second* one first so it
#     naturally inherits bullet 2's level-0 formatting; then insert the sub-level bullet
#     between them and bump only its level) ---
$para2 = $tr.Paragraphs(2)
$para2.InsertAfter("`rUnit test for all database functions")

$para2 = $tr.Paragraphs(2)
$para2.InsertAfter("`rBecause of parse tree complexity, Black Box tests cover use cases we wouldn" + [char]8217 + "t come up with")
$para3 = $tr.Paragraphs(3)
$para3.IndentLevel = 2

# --- Final bullet: split "The tests use a common file..." into 3 runs ---
$para5 = $tr.Paragraphs(5)
$sub = $para5.Characters(5, 6)
$sub.Text = "Black Box tests "
